$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Direct estimate of R")

# The "Reference(s)" values for the Savoca et al. rows had been typed into
# column C - move them (value + style) over to column D where they belong,
# and remove the now-empty source cells.
$ws1.Range("C8").Cut($ws1.Range("D8"))
$ws1.Range("C8").Clear()

$ws1.Range("C9").Cut($ws1.Range("D9"))
$ws1.Range("C9").Clear()

$ws1.Range("C10").Cut($ws1.Range("D10"))
$ws1.Range("C10").Clear()

# Update the active selection to reflect where editing left off.
$ws1.Range("D15").Select()
